# Apply the updated cryptocurrency market data pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '74.802.04'
$ws.Range('E2').Value = '  -0.05%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.844.82'
$ws.Range('E3').Value = '  +9.48%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '189.75'
$ws.Range('E5').Value = '  +1.79%  '

# Row 6
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '602.18'
$ws.Range('E6').Value = '  +3.35%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.560'
$ws.Range('E8').Value = '  +4.72%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.194'
$ws.Range('E9').Value = '  -7.88%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.842.64'
$ws.Range('E10').Value = '  +9.32%  '

# Row 11
$ws.Range('E11').Value = '  +0.00%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.373'
$ws.Range('E12').Value = '  +3.17%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.93'
$ws.Range('E13').Value = '  +2.15%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.363.94'
$ws.Range('E14').Value = '  +10.01%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '74.856.37'
$ws.Range('E15').Value = '  +0.34%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.46'
$ws.Range('E16').Value = '  +4.13%  '

# Row 17
$ws.Range('E17').Value = '  -2.48%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.849.33'
$ws.Range('E18').Value = '  +9.80%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.20'
$ws.Range('E19').Value = '  +8.20%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.46'
$ws.Range('E20').Value = '  +5.98%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '376.83'
$ws.Range('E21').Value = '  -0.45%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.28'
$ws.Range('E22').Value = '  -0.68%  '

# Row 23
$ws.Range('E23').Value = '  +1.45%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.21'
$ws.Range('E24').Value = '  -0.50%  '

# Row 25
$ws.Range('E25').Value = '  +0.00%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.85'
$ws.Range('E26').Value = '  +0.91%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.24'
$ws.Range('E27').Value = '  +0.53%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.991.36'
$ws.Range('E28').Value = '  +9.55%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.65'
$ws.Range('E29').Value = '  +4.53%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000104'
$ws.Range('E30').Value = '  +10.16%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.04%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '529.27'
$ws.Range('E32').Value = '  +4.93%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.40'
$ws.Range('E33').Value = '  +4.69%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.93'
$ws.Range('E34').Value = '  -0.63%  '

# Row 35
$ws.Range('E35').Value = '  +5.88%  '

# Row 36
$ws.Range('E36').Value = '  -0.11%  '

# Row 37
$ws.Range('E37').Value = '  +1.02%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '20.14'
$ws.Range('E38').Value = '  +4.53%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '162.15'
$ws.Range('E39').Value = '  +1.28%  '

# Row 40
$ws.Range('E40').Value = '  -0.63%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '182.62'
$ws.Range('E41').Value = '  +22.48%  '

# Row 42
$ws.Range('E42').Value = '  +0.01%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.12'
$ws.Range('E43').Value = '  +2.02%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.341'
$ws.Range('E44').Value = '  +6.10%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.70'
$ws.Range('E45').Value = '  +0.72%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.25'
$ws.Range('E46').Value = '  +7.63%  '

# Row 47
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.71'
$ws.Range('E47').Value = '  +1.57%  '

# Row 48
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.38'
$ws.Range('E48').Value = '  -2.58%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0851'
$ws.Range('E49').Value = '  +4.38%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.574'
$ws.Range('E50').Value = '  +9.85%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.77'
$ws.Range('E51').Value = '  +4.02%  '
